# Auto-generated edit script applying numeric corrections to multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 947.5  # H38: 903.41174 -> 947.5
$ws.Cells.Item(38, 9).Value = 947.5  # I38: 1010.5333 -> 947.5
$ws.Cells.Item(38, 10).Value = 0  # J38: 100 -> 0
$ws.Cells.Item(38, 11).Value = 2842.5  # K38: 3031.5999 -> 2842.5
$ws.Cells.Item(38, 12).Value = 0  # L38: 300 -> 0
$ws.Cells.Item(38, 13).Value = -2470.5  # M38: -2659.5999 -> -2470.5
$ws.Cells.Item(38, 14).ClearContents()  # N38: -1044 -> (removed)
$ws.Cells.Item(51, 8).Value = 18612.416  # H51: 26609.674 -> 18612.416
$ws.Cells.Item(51, 10).Value = 20587.889  # J51: 41216.25 -> 20587.889
$ws.Cells.Item(51, 12).Value = 20587.889  # L51: 41216.25 -> 20587.889
$ws.Cells.Item(51, 14).Value = -21555.889  # N51: -42184.25 -> -21555.889
$ws.Cells.Item(62, 8).Value = 5024  # H62: 5141.857 -> 5024
$ws.Cells.Item(62, 10).Value = 4247.25  # J62: 4263.3335 -> 4247.25
$ws.Cells.Item(62, 12).Value = 4247.25  # L62: 4263.3335 -> 4247.25
$ws.Cells.Item(62, 14).Value = -5495.25  # N62: -5511.3335 -> -5495.25
$ws.Cells.Item(65, 8).Value = 5024  # H65: 5141.857 -> 5024
$ws.Cells.Item(65, 10).Value = 4247.25  # J65: 4263.3335 -> 4247.25
$ws.Cells.Item(65, 12).Value = 21236.25  # L65: 21316.6675 -> 21236.25
$ws.Cells.Item(65, 14).Value = -27476.25  # N65: -27556.6675 -> -27476.25
$ws.Cells.Item(76, 8).Value = 5213.5  # H76: 5070.8 -> 5213.5
$ws.Cells.Item(76, 10).Value = 4977  # J76: 4818 -> 4977
$ws.Cells.Item(76, 12).Value = 4977  # L76: 4818 -> 4977
$ws.Cells.Item(76, 14).Value = -5607  # N76: -5448 -> -5607
$ws.Cells.Item(79, 8).Value = 5213.5  # H79: 5070.8 -> 5213.5
$ws.Cells.Item(79, 10).Value = 4977  # J79: 4818 -> 4977
$ws.Cells.Item(79, 12).Value = 4977  # L79: 4818 -> 4977
$ws.Cells.Item(79, 14).Value = -7161  # N79: -7002 -> -7161
$ws.Cells.Item(112, 8).Value = 1729.5  # H112: 1895.4 -> 1729.5
$ws.Cells.Item(112, 10).Value = 1729.5  # J112: 1895.4 -> 1729.5
$ws.Cells.Item(112, 12).Value = 5188.5  # L112: 5686.200000000001 -> 5188.5
$ws.Cells.Item(112, 14).Value = -7404.5  # N112: -7902.200000000001 -> -7404.5
$ws.Cells.Item(129, 8).Value = 1835.4773  # H129: 1782.1224 -> 1835.4773
$ws.Cells.Item(129, 9).Value = 1463  # I129: 1350.2858 -> 1463
$ws.Cells.Item(129, 10).Value = 1862.7317  # J129: 1854.0952 -> 1862.7317
$ws.Cells.Item(129, 11).Value = 4389  # K129: 4050.8574 -> 4389
$ws.Cells.Item(129, 12).Value = 5588.1951  # L129: 5562.2856 -> 5588.1951
$ws.Cells.Item(129, 13).Value = 611  # M129: 949.1425999999997 -> 611
$ws.Cells.Item(129, 14).Value = -15588.1951  # N129: -15562.2856 -> -15588.1951
$ws.Cells.Item(132, 8).Value = 40977.08  # H132: 44532 -> 40977.08
$ws.Cells.Item(132, 9).Value = 44447.824  # I132: 48671.855 -> 44447.824
$ws.Cells.Item(132, 11).Value = 133343.472  # K132: 146015.565 -> 133343.472
$ws.Cells.Item(132, 13).Value = -130813.472  # M132: -143485.565 -> -130813.472
$ws.Cells.Item(137, 9).Value = 3434  # I137: 3720.8 -> 3434
$ws.Cells.Item(137, 10).Value = 10667.667  # J137: 8500.75 -> 10667.667
$ws.Cells.Item(137, 11).Value = 10302  # K137: 11162.4 -> 10302
$ws.Cells.Item(137, 12).Value = 32003.001  # L137: 25502.25 -> 32003.001
$ws.Cells.Item(137, 13).Value = -7752  # M137: -8612.400000000001 -> -7752
$ws.Cells.Item(137, 14).Value = -37103.001  # N137: -30602.25 -> -37103.001
$ws.Cells.Item(138, 8).Value = 2354.0557  # H138: 2236.0625 -> 2354.0557
$ws.Cells.Item(138, 9).Value = 1971.6154  # I138: 1979.1538 -> 1971.6154
$ws.Cells.Item(138, 10).Value = 3348.4  # J138: 3349.3333 -> 3348.4
$ws.Cells.Item(138, 11).Value = 5914.8462  # K138: 5937.4614 -> 5914.8462
$ws.Cells.Item(138, 12).Value = 10045.2  # L138: 10047.9999 -> 10045.2
$ws.Cells.Item(138, 13).Value = -774.8462  # M138: -797.4614000000001 -> -774.8462
$ws.Cells.Item(138, 14).Value = -20325.2  # N138: -20327.9999 -> -20325.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2340.6562  # H32: 2723.36 -> 2340.6562
$ws.Cells.Item(32, 9).Value = 2351.6453  # I32: 2753.5 -> 2351.6453
$ws.Cells.Item(32, 11).Value = 2351.6453  # K32: 2753.5 -> 2351.6453
$ws.Cells.Item(32, 13).Value = -2064.6453  # M32: -2466.5 -> -2064.6453
$ws.Cells.Item(61, 8).Value = 3434.238  # H61: 2749.2424 -> 3434.238
$ws.Cells.Item(61, 9).Value = 2435.3635  # I61: 2077.0952 -> 2435.3635
$ws.Cells.Item(61, 10).Value = 4533  # J61: 3925.5 -> 4533
$ws.Cells.Item(61, 11).Value = 2435.3635  # K61: 2077.0952 -> 2435.3635
$ws.Cells.Item(61, 12).Value = 4533  # L61: 3925.5 -> 4533
$ws.Cells.Item(61, 13).Value = -2223.3635  # M61: -1865.0952 -> -2223.3635
$ws.Cells.Item(61, 14).Value = -4957  # N61: -4349.5 -> -4957
$ws.Cells.Item(74, 8).Value = 7717554.5  # H74: 3087571.8 -> 7717554.5
$ws.Cells.Item(74, 9).Value = 3705510  # I74: 1544468 -> 3705510
$ws.Cells.Item(74, 10).Value = 27777776  # J74: 9259987 -> 27777776
$ws.Cells.Item(74, 11).Value = 3705510  # K74: 1544468 -> 3705510
$ws.Cells.Item(74, 12).Value = 27777776  # L74: 9259987 -> 27777776
$ws.Cells.Item(74, 13).Value = -3704636  # M74: -1543594 -> -3704636
$ws.Cells.Item(74, 14).Value = -27779524  # N74: -9261735 -> -27779524
$ws.Cells.Item(77, 8).Value = 7717554.5  # H77: 3087571.8 -> 7717554.5
$ws.Cells.Item(77, 9).Value = 3705510  # I77: 1544468 -> 3705510
$ws.Cells.Item(77, 10).Value = 27777776  # J77: 9259987 -> 27777776
$ws.Cells.Item(77, 11).Value = 18527550  # K77: 7722340 -> 18527550
$ws.Cells.Item(77, 12).Value = 138888880  # L77: 46299935 -> 138888880
$ws.Cells.Item(77, 13).Value = -18523182  # M77: -7717972 -> -18523182
$ws.Cells.Item(77, 14).Value = -138897616  # N77: -46308671 -> -138897616
$ws.Cells.Item(102, 8).Value = 624.43475  # H102: 618 -> 624.43475
$ws.Cells.Item(102, 9).Value = 488.3  # I102: 480.9 -> 488.3
$ws.Cells.Item(102, 11).Value = 488.3  # K102: 480.9 -> 488.3
$ws.Cells.Item(102, 13).Value = 1133.7  # M102: 1141.1 -> 1133.7
$ws.Cells.Item(112, 8).Value = 0  # H112: 69994 -> 0
$ws.Cells.Item(112, 10).Value = 0  # J112: 69994 -> 0
$ws.Cells.Item(112, 12).Value = 0  # L112: 69994 -> 0
$ws.Cells.Item(112, 14).ClearContents()  # N112: -72948 -> (removed)
$ws.Cells.Item(132, 8).Value = 17243598  # H132: 21741456 -> 17243598
$ws.Cells.Item(132, 9).Value = 2224.318  # I132: 2370.0557 -> 2224.318
$ws.Cells.Item(132, 10).Value = 71430776  # J132: 100002160 -> 71430776
$ws.Cells.Item(132, 11).Value = 6672.954000000001  # K132: 7110.1671 -> 6672.954000000001
$ws.Cells.Item(132, 12).Value = 214292328  # L132: 300006480 -> 214292328
$ws.Cells.Item(132, 13).Value = -4142.954000000001  # M132: -4580.1671 -> -4142.954000000001
$ws.Cells.Item(132, 14).Value = -214297388  # N132: -300011540 -> -214297388
$ws.Cells.Item(136, 8).Value = 3434.238  # H136: 2749.2424 -> 3434.238
$ws.Cells.Item(136, 9).Value = 2435.3635  # I136: 2077.0952 -> 2435.3635
$ws.Cells.Item(136, 10).Value = 4533  # J136: 3925.5 -> 4533
$ws.Cells.Item(136, 11).Value = 7306.0905  # K136: 6231.285600000001 -> 7306.0905
$ws.Cells.Item(136, 12).Value = 13599  # L136: 11776.5 -> 13599
$ws.Cells.Item(136, 13).Value = -4756.0905  # M136: -3681.285600000001 -> -4756.0905
$ws.Cells.Item(136, 14).Value = -18699  # N136: -16876.5 -> -18699
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 971.5714  # H20: 961.2 -> 971.5714
$ws.Cells.Item(20, 9).Value = 950.1667  # I20: 961.2 -> 950.1667
$ws.Cells.Item(20, 10).Value = 1100  # J20: 0 -> 1100
$ws.Cells.Item(20, 11).Value = 950.1667  # K20: 961.2 -> 950.1667
$ws.Cells.Item(20, 12).Value = 1100  # L20: 0 -> 1100
$ws.Cells.Item(20, 13).Value = -703.1667  # M20: -714.2 -> -703.1667
$ws.Cells.Item(20, 14).Value = -1594  # N20: None -> -1594
$ws.Cells.Item(35, 8).Value = 41995  # H35: 42000 -> 41995
$ws.Cells.Item(35, 10).Value = 41995  # J35: 42000 -> 41995
$ws.Cells.Item(35, 12).Value = 41995  # L35: 42000 -> 41995
$ws.Cells.Item(35, 14).Value = -42615  # N35: -42620 -> -42615
$ws.Cells.Item(86, 8).Value = 3116.7  # H86: 3006.818 -> 3116.7
$ws.Cells.Item(86, 9).Value = 2458  # I86: 2318.5 -> 2458
$ws.Cells.Item(86, 10).Value = 3399  # J86: 3400.1428 -> 3399
$ws.Cells.Item(86, 11).Value = 2458  # K86: 2318.5 -> 2458
$ws.Cells.Item(86, 12).Value = 3399  # L86: 3400.1428 -> 3399
$ws.Cells.Item(86, 13).Value = -1335  # M86: -1195.5 -> -1335
$ws.Cells.Item(86, 14).Value = -5645  # N86: -5646.1428 -> -5645
$ws.Cells.Item(89, 8).Value = 3116.7  # H89: 3006.818 -> 3116.7
$ws.Cells.Item(89, 9).Value = 2458  # I89: 2318.5 -> 2458
$ws.Cells.Item(89, 10).Value = 3399  # J89: 3400.1428 -> 3399
$ws.Cells.Item(89, 11).Value = 12290  # K89: 11592.5 -> 12290
$ws.Cells.Item(89, 12).Value = 16995  # L89: 17000.714 -> 16995
$ws.Cells.Item(89, 13).Value = -6674  # M89: -5976.5 -> -6674
$ws.Cells.Item(89, 14).Value = -28227  # N89: -28232.714 -> -28227
$ws.Cells.Item(99, 8).Value = 2624.12  # H99: 2929.3044 -> 2624.12
$ws.Cells.Item(99, 9).Value = 2534.3635  # I99: 2894.3 -> 2534.3635
$ws.Cells.Item(99, 10).Value = 3282.3333  # J99: 3162.6667 -> 3282.3333
$ws.Cells.Item(99, 11).Value = 2534.3635  # K99: 2894.3 -> 2534.3635
$ws.Cells.Item(99, 12).Value = 3282.3333  # L99: 3162.6667 -> 3282.3333
$ws.Cells.Item(99, 13).Value = -1036.3635  # M99: -1396.3 -> -1036.3635
$ws.Cells.Item(99, 14).Value = -6278.3333  # N99: -6158.6667 -> -6278.3333
$ws.Cells.Item(134, 8).Value = 25364628  # H134: 32410074 -> 25364628
$ws.Cells.Item(134, 9).Value = 13160162  # I134: 17859822 -> 13160162
$ws.Cells.Item(134, 10).Value = 83335830  # J134: 83335960 -> 83335830
$ws.Cells.Item(134, 11).Value = 39480486  # K134: 53579466 -> 39480486
$ws.Cells.Item(134, 12).Value = 250007490  # L134: 250007880 -> 250007490
$ws.Cells.Item(134, 13).Value = -39477951  # M134: -53576931 -> -39477951
$ws.Cells.Item(134, 14).Value = -250012560  # N134: -250012950 -> -250012560
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2392.5454  # H31: 1925.9375 -> 2392.5454
$ws.Cells.Item(31, 9).Value = 2211.8  # I31: 1774.3334 -> 2211.8
$ws.Cells.Item(31, 11).Value = 2211.8  # K31: 1774.3334 -> 2211.8
$ws.Cells.Item(31, 13).Value = -1916.8  # M31: -1479.3334 -> -1916.8
$ws.Cells.Item(34, 8).Value = 2392.5454  # H34: 1925.9375 -> 2392.5454
$ws.Cells.Item(34, 9).Value = 2211.8  # I34: 1774.3334 -> 2211.8
$ws.Cells.Item(34, 11).Value = 2211.8  # K34: 1774.3334 -> 2211.8
$ws.Cells.Item(34, 13).Value = -2009.8  # M34: -1572.3334 -> -2009.8
$ws.Cells.Item(86, 8).Value = 3151.5715  # H86: 3162.25 -> 3151.5715
$ws.Cells.Item(86, 9).Value = 3010.1667  # I86: 2883 -> 3010.1667
$ws.Cells.Item(86, 11).Value = 3010.1667  # K86: 2883 -> 3010.1667
$ws.Cells.Item(86, 13).Value = -1887.1667  # M86: -1760 -> -1887.1667
$ws.Cells.Item(89, 8).Value = 3151.5715  # H89: 3162.25 -> 3151.5715
$ws.Cells.Item(89, 9).Value = 3010.1667  # I89: 2883 -> 3010.1667
$ws.Cells.Item(89, 11).Value = 15050.8335  # K89: 14415 -> 15050.8335
$ws.Cells.Item(89, 13).Value = -9434.833500000001  # M89: -8799 -> -9434.833500000001
$ws.Cells.Item(134, 8).Value = 5559015  # H134: 6253788 -> 5559015
$ws.Cells.Item(134, 9).Value = 2024.7  # I134: 2323.25 -> 2024.7
$ws.Cells.Item(134, 11).Value = 6074.1  # K134: 6969.75 -> 6074.1
$ws.Cells.Item(134, 13).Value = -3539.1  # M134: -4434.75 -> -3539.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 763.5714  # H5: 669 -> 763.5714
$ws.Cells.Item(5, 9).Value = 1140.5555  # I5: 878 -> 1140.5555
$ws.Cells.Item(5, 10).Value = 480.83334  # J5: 397.3 -> 480.83334
$ws.Cells.Item(5, 11).Value = 3421.6665  # K5: 2634 -> 3421.6665
$ws.Cells.Item(5, 12).Value = 1442.50002  # L5: 1191.9 -> 1442.50002
$ws.Cells.Item(5, 13).Value = -3309.6665  # M5: -2522 -> -3309.6665
$ws.Cells.Item(5, 14).Value = -1666.50002  # N5: -1415.9 -> -1666.50002
$ws.Cells.Item(107, 8).Value = 843.8333  # H107: 749.7059 -> 843.8333
$ws.Cells.Item(107, 9).Value = 374.2  # I107: 377.2 -> 374.2
$ws.Cells.Item(107, 10).Value = 1179.2858  # J107: 904.9167 -> 1179.2858
$ws.Cells.Item(107, 11).Value = 1122.6  # K107: 1131.6 -> 1122.6
$ws.Cells.Item(107, 12).Value = 3537.8574  # L107: 2714.7501 -> 3537.8574
$ws.Cells.Item(107, 13).Value = 797.4000000000001  # M107: 788.4000000000001 -> 797.4000000000001
$ws.Cells.Item(107, 14).Value = -7377.857400000001  # N107: -6554.7501 -> -7377.857400000001
$ws.Cells.Item(117, 8).Value = 2932.7334  # H117: 3192.1538 -> 2932.7334
$ws.Cells.Item(117, 9).Value = 1999.6666  # I117: 1749.75 -> 1999.6666
$ws.Cells.Item(117, 10).Value = 3166  # J117: 3833.2222 -> 3166
$ws.Cells.Item(117, 11).Value = 5998.9998  # K117: 5249.25 -> 5998.9998
$ws.Cells.Item(117, 12).Value = 9498  # L117: 11499.6666 -> 9498
$ws.Cells.Item(117, 13).Value = -2556.9998  # M117: -1807.25 -> -2556.9998
$ws.Cells.Item(117, 14).Value = -16382  # N117: -18383.6666 -> -16382
$ws.Cells.Item(135, 8).Value = 763.5714  # H135: 669 -> 763.5714
$ws.Cells.Item(135, 9).Value = 1140.5555  # I135: 878 -> 1140.5555
$ws.Cells.Item(135, 10).Value = 480.83334  # J135: 397.3 -> 480.83334
$ws.Cells.Item(135, 11).Value = 10264.9995  # K135: 7902 -> 10264.9995
$ws.Cells.Item(135, 12).Value = 4327.50006  # L135: 3575.7 -> 4327.50006
$ws.Cells.Item(135, 13).Value = -7729.9995  # M135: -5367 -> -7729.9995
$ws.Cells.Item(135, 14).Value = -9397.50006  # N135: -8645.700000000001 -> -9397.50006
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 10042  # H26: 16361.333 -> 10042
$ws.Cells.Item(26, 10).Value = 10042  # J26: 16361.333 -> 10042
$ws.Cells.Item(26, 12).Value = 10042  # L26: 16361.333 -> 10042
$ws.Cells.Item(26, 14).Value = -10602  # N26: -16921.333 -> -10602
$ws.Cells.Item(50, 8).Value = 10042  # H50: 16361.333 -> 10042
$ws.Cells.Item(50, 10).Value = 10042  # J50: 16361.333 -> 10042
$ws.Cells.Item(50, 12).Value = 10042  # L50: 16361.333 -> 10042
$ws.Cells.Item(50, 14).Value = -11038  # N50: -17357.333 -> -11038
$ws.Cells.Item(80, 8).Value = 5381.364  # H80: 7391.3335 -> 5381.364
$ws.Cells.Item(80, 9).Value = 3419.5  # I80: 3869.6 -> 3419.5
$ws.Cells.Item(80, 11).Value = 3419.5  # K80: 3869.6 -> 3419.5
$ws.Cells.Item(80, 13).Value = -2421.5  # M80: -2871.6 -> -2421.5
$ws.Cells.Item(83, 8).Value = 5381.364  # H83: 7391.3335 -> 5381.364
$ws.Cells.Item(83, 9).Value = 3419.5  # I83: 3869.6 -> 3419.5
$ws.Cells.Item(83, 11).Value = 17097.5  # K83: 19348 -> 17097.5
$ws.Cells.Item(83, 13).Value = -12105.5  # M83: -14356 -> -12105.5
$ws.Cells.Item(102, 8).Value = 2232.5  # H102: 2264.5833 -> 2232.5
$ws.Cells.Item(102, 9).Value = 2271.818  # I102: 2306.818 -> 2271.818
$ws.Cells.Item(102, 11).Value = 2271.818  # K102: 2306.818 -> 2271.818
$ws.Cells.Item(102, 13).Value = -649.8180000000002  # M102: -684.8180000000002 -> -649.8180000000002
$ws.Cells.Item(132, 8).Value = 1942  # H132: 1759.1538 -> 1942
$ws.Cells.Item(132, 9).Value = 1909.75  # I132: 1715.5454 -> 1909.75
$ws.Cells.Item(132, 10).Value = 2200  # J132: 1999 -> 2200
$ws.Cells.Item(132, 11).Value = 5729.25  # K132: 5146.6362 -> 5729.25
$ws.Cells.Item(132, 12).Value = 6600  # L132: 5997 -> 6600
$ws.Cells.Item(132, 13).Value = -3199.25  # M132: -2616.6362 -> -3199.25
$ws.Cells.Item(132, 14).Value = -11660  # N132: -11057 -> -11660
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1116.4865  # H22: 1233.4375 -> 1116.4865
$ws.Cells.Item(22, 9).Value = 984.3  # I22: 1058.5294 -> 984.3
$ws.Cells.Item(22, 10).Value = 1272  # J22: 1431.6666 -> 1272
$ws.Cells.Item(22, 11).Value = 984.3  # K22: 1058.5294 -> 984.3
$ws.Cells.Item(22, 12).Value = 1272  # L22: 1431.6666 -> 1272
$ws.Cells.Item(22, 13).Value = -689.3  # M22: -763.5293999999999 -> -689.3
$ws.Cells.Item(22, 14).Value = -1862  # N22: -2021.6666 -> -1862
$ws.Cells.Item(27, 8).Value = 1116.4865  # H27: 1233.4375 -> 1116.4865
$ws.Cells.Item(27, 9).Value = 984.3  # I27: 1058.5294 -> 984.3
$ws.Cells.Item(27, 10).Value = 1272  # J27: 1431.6666 -> 1272
$ws.Cells.Item(27, 11).Value = 984.3  # K27: 1058.5294 -> 984.3
$ws.Cells.Item(27, 12).Value = 1272  # L27: 1431.6666 -> 1272
$ws.Cells.Item(27, 13).Value = -877.3  # M27: -951.5293999999999 -> -877.3
$ws.Cells.Item(27, 14).Value = -1486  # N27: -1645.6666 -> -1486
$ws.Cells.Item(46, 8).Value = 2652.6562  # H46: 2593.4849 -> 2652.6562
$ws.Cells.Item(46, 9).Value = 1766.6666  # I46: 1500 -> 1766.6666
$ws.Cells.Item(46, 11).Value = 1766.6666  # K46: 1500 -> 1766.6666
$ws.Cells.Item(46, 13).Value = -1578.6666  # M46: -1312 -> -1578.6666
$ws.Cells.Item(97, 8).Value = 11109.6  # H97: 10141.333 -> 11109.6
$ws.Cells.Item(97, 10).Value = 11109.6  # J97: 10141.333 -> 11109.6
$ws.Cells.Item(97, 12).Value = 11109.6  # L97: 10141.333 -> 11109.6
$ws.Cells.Item(97, 14).Value = -13091.6  # N97: -12123.333 -> -13091.6
$ws.Cells.Item(122, 8).Value = 3568.5833  # H122: 3280.5454 -> 3568.5833
$ws.Cells.Item(122, 9).Value = 3315.75  # I122: 3000.25 -> 3315.75
$ws.Cells.Item(122, 10).Value = 3695  # J122: 3616.9 -> 3695
$ws.Cells.Item(122, 11).Value = 9947.25  # K122: 9000.75 -> 9947.25
$ws.Cells.Item(122, 12).Value = 11085  # L122: 10850.7 -> 11085
$ws.Cells.Item(122, 13).Value = -7497.25  # M122: -6550.75 -> -7497.25
$ws.Cells.Item(122, 14).Value = -15985  # N122: -15750.7 -> -15985
$ws.Cells.Item(132, 8).Value = 578.6  # H132: 0 -> 578.6
$ws.Cells.Item(132, 9).Value = 631.3333  # I132: 0 -> 631.3333
$ws.Cells.Item(132, 10).Value = 499.5  # J132: 0 -> 499.5
$ws.Cells.Item(132, 11).Value = 1893.9999  # K132: 0 -> 1893.9999
$ws.Cells.Item(132, 12).Value = 1498.5  # L132: 0 -> 1498.5
$ws.Cells.Item(132, 13).Value = 636.0001  # M132: None -> 636.0001
$ws.Cells.Item(132, 14).Value = -6558.5  # N132: None -> -6558.5
$ws.Cells.Item(136, 8).Value = 76927420  # H136: 71432344 -> 76927420
$ws.Cells.Item(136, 9).Value = 4773.2  # I136: 3222.8572 -> 4773.2
$ws.Cells.Item(136, 10).Value = 125004080  # J136: 142861470 -> 125004080
$ws.Cells.Item(136, 11).Value = 14319.6  # K136: 9668.571599999999 -> 14319.6
$ws.Cells.Item(136, 12).Value = 375012240  # L136: 428584410 -> 375012240
$ws.Cells.Item(136, 13).Value = -11769.6  # M136: -7118.571599999999 -> -11769.6
$ws.Cells.Item(136, 14).Value = -375017340  # N136: -428589510 -> -375017340
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(116, 8).Value = 42500  # H116: 70000 -> 42500
$ws.Cells.Item(116, 10).Value = 42500  # J116: 70000 -> 42500
$ws.Cells.Item(116, 12).Value = 42500  # L116: 70000 -> 42500
$ws.Cells.Item(116, 14).Value = -51678  # N116: -79178 -> -51678
$ws.Cells.Item(126, 8).Value = 2192.1538  # H126: 2273.0908 -> 2192.1538
$ws.Cells.Item(126, 9).Value = 1328.1666  # I126: 1395.2 -> 1328.1666
$ws.Cells.Item(126, 10).Value = 2932.7144  # J126: 3004.6667 -> 2932.7144
$ws.Cells.Item(126, 11).Value = 3984.4998  # K126: 4185.6 -> 3984.4998
$ws.Cells.Item(126, 12).Value = 8798.143199999999  # L126: 9014.000100000001 -> 8798.143199999999
$ws.Cells.Item(126, 13).Value = -1514.4998  # M126: -1715.6 -> -1514.4998
$ws.Cells.Item(126, 14).Value = -13738.1432  # N126: -13954.0001 -> -13738.1432
